$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D5").Value = "1종 오류와 2종 오류"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2021/01/26/types_of_errors.html"

$ws.Range("D16").Value = "주요 기계학습, 딥러닝 activation function"
$ws.Range("E16").Value = "https://wewinserv.tistory.com/141"

$ws.Range("D22").Value = "Spectral GCN 은… 사드세요"
$ws.Range("E22").Value = "https://tootouch.github.io/research/spectral_gcn/"

$ws.Range("D37").Value = "[Paper Review] Transferable Multi-Domain State Generator for Task-Oriented Dialogue Systems"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1433&mod=document&pageid=1"

$ws.Range("D39").Value = "Deep Face Detection with OpenCV in Python"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Deep-Face-Detection-with-OpenCV-in-Python-1"

$wb.Save()
